# Auto-generated edit script applying the row-level data corrections
# described by the diff for rows 11-23 (species-occurrence records).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Range("A11").Value = 111577080
$ws.Range("B11").Value = 89405
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 1202
$ws.Range("F11").Value = "Ullticka"
$ws.Range("G11").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H11").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q11").Value = 562937.8525077751
$ws.Range("R11").Value = 6954467.524316943
$ws.Range("Z11").Value = "17:24"
$ws.Range("AB11").Value = "17:24"

# Row 12
$ws.Range("A12").Value = 111577804
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = "VU"
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = "Knärot"
$ws.Range("G12").Value = "Goodyera repens"
$ws.Range("H12").Value = "(L.) R. Br."
$ws.Range("Q12").Value = 562814.775380839
$ws.Range("R12").Value = 6954390.834027934

# Row 13
$ws.Range("A13").Value = 111577180
$ws.Range("B13").Value = 89686
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 658
$ws.Range("F13").Value = "Rosenticka"
$ws.Range("G13").Value = "Rhodofomes roseus"
$ws.Range("H13").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q13").Value = 562883.198926247
$ws.Range("R13").Value = 6954441.700568204
$ws.Range("Z13").Value = "17:24"
$ws.Range("AB13").Value = "17:24"

# Row 14
$ws.Range("A14").Value = 111577347
$ws.Range("B14").Value = 96348
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("Q14").Value = 562796.503171768
$ws.Range("R14").Value = 6954336.792844097
$ws.Range("Z14").Value = "00:00"
$ws.Range("AB14").Value = "00:00"

# Row 15
$ws.Range("A15").Value = 111577964
$ws.Range("B15").Value = 96348
$ws.Range("D15").Value = "VU"
$ws.Range("E15").Value = 220787
$ws.Range("F15").Value = "Knärot"
$ws.Range("G15").Value = "Goodyera repens"
$ws.Range("H15").Value = "(L.) R. Br."
$ws.Range("Q15").Value = 562836.2055113926
$ws.Range("R15").Value = 6954423.824987715
$ws.Range("Z15").Value = "18:23"
$ws.Range("AB15").Value = "18:23"
$ws.Range("M15").ClearContents()

# Row 16
$ws.Range("A16").Value = 111577838
$ws.Range("Q16").Value = 562801.9298062191
$ws.Range("R16").Value = 6954389.67147268
$ws.Range("Z16").Value = "18:14"
$ws.Range("AB16").Value = "18:14"

# Row 17
$ws.Range("A17").Value = 111575408
$ws.Range("B17").Value = 96348
$ws.Range("D17").Value = "VU"
$ws.Range("E17").Value = 220787
$ws.Range("F17").Value = "Knärot"
$ws.Range("G17").Value = "Goodyera repens"
$ws.Range("H17").Value = "(L.) R. Br."
$ws.Range("Q17").Value = 562539.034657649
$ws.Range("R17").Value = 6954609.073577877
$ws.Range("Z17").Value = "16:21"
$ws.Range("AB17").Value = "16:21"

# Row 18
$ws.Range("A18").Value = 111577591
$ws.Range("B18").Value = 56543
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 103021
$ws.Range("F18").Value = "Talltita"
$ws.Range("G18").Value = "Poecile montanus"
$ws.Range("H18").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("Q18").Value = 562822.1033927511
$ws.Range("R18").Value = 6954368.028004575
$ws.Range("Z18").Value = "00:00"
$ws.Range("AB18").Value = "00:00"
$ws.Range("M18").Value = "lockläte, övriga läten"

# Row 19
$ws.Range("A19").Value = 111577611
$ws.Range("B19").Value = 96348
$ws.Range("E19").Value = 220787
$ws.Range("F19").Value = "Knärot"
$ws.Range("G19").Value = "Goodyera repens"
$ws.Range("H19").Value = "(L.) R. Br."
$ws.Range("Q19").Value = 562810.9079359611
$ws.Range("R19").Value = 6954400.856378952
$ws.Range("Z19").Value = "00:00"
$ws.Range("AB19").Value = "00:00"

# Row 20
$ws.Range("A20").Value = 111577273
$ws.Range("B20").Value = 89686
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 658
$ws.Range("F20").Value = "Rosenticka"
$ws.Range("G20").Value = "Rhodofomes roseus"
$ws.Range("H20").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q20").Value = 562824.7977144517
$ws.Range("R20").Value = 6954323.105396069
$ws.Range("Z20").Value = "17:24"
$ws.Range("AB20").Value = "17:24"

# Row 21
$ws.Range("A21").Value = 111577743
$ws.Range("Q21").Value = 562802.8660743404
$ws.Range("R21").Value = 6954388.771485241
$ws.Range("Z21").Value = "17:59"
$ws.Range("AB21").Value = "17:59"

# Row 22
$ws.Range("A22").Value = 111577193
$ws.Range("B22").Value = 89845
$ws.Range("D22").Value = "VU"
$ws.Range("E22").Value = 1209
$ws.Range("F22").Value = "Rynkskinn"
$ws.Range("G22").Value = "Phlebia centrifuga"
$ws.Range("H22").Value = "P.Karst."
$ws.Range("Q22").Value = 562885.3077477051
$ws.Range("R22").Value = 6954427.514711756

# Row 23
$ws.Range("A23").Value = 111577994
$ws.Range("B23").Value = 89686
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 658
$ws.Range("F23").Value = "Rosenticka"
$ws.Range("G23").Value = "Rhodofomes roseus"
$ws.Range("H23").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("Q23").Value = 562890.6102569005
$ws.Range("R23").Value = 6954486.814324431
$ws.Range("Z23").Value = "00:00"
$ws.Range("AB23").Value = "00:00"
